$d = $word.ActiveDocument

# Styling bug: the "Commandtax Examples" source-code paragraph carried
# literal Markdown code-fence markers ("```") as their own runs (each
# on its own line, via a manual line break) inside the Word doc. Strip
# the opening "```" + following line break and the trailing line break
# + closing "```" from that paragraph, leaving the real command lines
# (and the line breaks between them) exactly as they were.

$vt = [char]11   # manual line break char, i.e. <w:br w:type="textWrapping"/>

# Find the SourceCode paragraph that starts and ends with the fence.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $trimmed = $t.TrimEnd([char]13)
    if ($t.StartsWith('```') -and $trimmed.EndsWith('```')) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Remove the opening fence: "```" immediately followed by a break.
    $rng = $target.Range
    $openFence = '```' + $vt
    $foundOpen = $rng.Find.Execute($openFence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($foundOpen) {
        $rng.Text = ""
    }

    # Remove the closing fence: a break immediately followed by "```".
    $rng2 = $target.Range
    $closeFence = $vt + '```'
    $foundClose = $rng2.Find.Execute($closeFence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($foundClose) {
        $rng2.Text = ""
    }

    Write-Output "openFenceRemoved=$foundOpen closeFenceRemoved=$foundClose"
} else {
    Write-Output "target paragraph not found"
}
